# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.021.03"
$ws.Range("E2").Value = "  +0.92%  "

$ws.Range("D3").Value = "2.297.90"
$ws.Range("E3").Value = "  +0.45%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.56"
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.52"
$ws.Range("E6").Value = "  +0.56%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.506"
$ws.Range("E7").Value = "  +0.42%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.505"
$ws.Range("E9").Value = "  +0.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.80"
$ws.Range("E10").Value = "  +1.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.04"
$ws.Range("E12").Value = "  -2.19%  "

$ws.Range("E13").Value = "  +2.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.08"
$ws.Range("E14").Value = "  +10.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.77"
$ws.Range("E15").Value = "  +1.73%  "

$ws.Range("D16").Value = "2.658.78"
$ws.Range("E16").Value = "  +0.64%  "

$ws.Range("D17").Value = "2.274.20"
$ws.Range("E17").Value = "  -1.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.810"
$ws.Range("E18").Value = "  +2.82%  "

$ws.Range("D19").Value = "42.933.99"
$ws.Range("E19").Value = "  +0.89%  "

$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.64"
$ws.Range("E20").Value = "  +0.95%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0903"
$ws.Range("E21").Value = "  +0.80%  "

$ws.Range("E22").Value = "  +0.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.50"
$ws.Range("E23").Value = "  +0.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.63"
$ws.Range("E24").Value = "  +0.78%  "

$ws.Range("E25").Value = "  +4.81%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E27").Value = "  -1.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.35"
$ws.Range("E28").Value = "  -0.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.71"
$ws.Range("E29").Value = "  +0.22%  "

$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.86"
$ws.Range("E31").Value = "  -0.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.12"
$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("E34").Value = "  -0.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.59"
$ws.Range("E35").Value = "  +5.67%  "

$ws.Range("E36").Value = "  +1.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.85"
$ws.Range("E37").Value = "  +3.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0702"
$ws.Range("E38").Value = "  +0.80%  "

$ws.Range("E39").Value = "  -0.39%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.78"
$ws.Range("E40").Value = "  +0.56%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.101"
$ws.Range("E41").Value = "  +0.46%  "

$ws.Range("E42").Value = "  -0.46%  "

$ws.Range("E43").Value = "  -4.36%  "

$ws.Range("D44").Value = "1.991.75"
$ws.Range("E44").Value = "  +1.39%  "

$ws.Range("E45").Value = "  +0.52%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.81"
$ws.Range("E46").Value = "  +1.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.48"
$ws.Range("E47").Value = "  -1.69%  "

$ws.Range("E48").Value = "  +0.58%  "

$ws.Range("D49").Value = "2.521.89"
$ws.Range("E49").Value = "  +0.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.33"
$ws.Range("E50").Value = "  +0.20%  "

$ws.Range("E51").Value = "  -1.84%  "
